# Add a "Media" (average) column to the BER simulation table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell V1: bold, yellow-filled, "0.00" number format, text "Media" ---
$header = $ws.Range("V1")
$header.Value = "Media"
$header.Font.Bold = $true
$header.Interior.Color = 65535
$header.NumberFormat = "0.00"

# --- Data cells V2:V31: per-row AVERAGE(B:U) formula, "0.00" number format ---
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("V$r")
    $cell.Formula = "=AVERAGE(B${r}:U${r})"
    $cell.NumberFormat = "0.00"
}

# --- Match the author's final selection: the whole column V highlighted ---
[void]$ws.Columns.Item(22).Select()
